$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A68").Value = "2021/08/16"
$ws.Range("B68").Value = 314.5
$ws.Range("C68").Value = 318.8
$ws.Range("D68").Value = 0.96
$ws.Range("E68").Value = 0.96
